# Add a new row for the "epitraxr" tool, inserted alphabetically before
# "epiworld: Fast Agent-Based Epi Models" (which currently sits at row 21),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21; everything from row 21 downward
# (epiworld, epiworld-forecasts, epiworldpy, ...) shifts down to make room.
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the epitraxr tool entry.
$ws.Range("A21").Value = "epitraxr: Manipulate Epitrax Data And Generate Reports"
$ws.Range("B21").Value = "R package for manipulating Epitrax data and generating reports. This tool is designed to simplify the process of working with Epitrax data, making it easier for public health professionals to analyze and report on disease surveillance data. It generates internal and public reports in CSV, Excel, or PDF formats."
$ws.Range("C21").Value = "Andrew Pulsipher"
$ws.Range("D21").Value = "a.pulsipher@utah.edu"
$ws.Range("E21").Value = "Yes"
$ws.Range("G21").Value = "Published"
$ws.Range("H21").Value = "MIT"
$ws.Range("I21").Value = "R"
$ws.Range("J21").Value = "Public Health Professionals"
$ws.Range("K21").Value = "Beginner"
$ws.Range("L21").Value = "Decision Support tool"
$ws.Range("M21").Value = "Epitrax data"
$ws.Range("N21").Value = "https://epiforesite.github.io/epitraxr/, https://github.com/EpiForeSITE/epitraxr"
$ws.Range("O21").Value = "https://github.com/EpiForeSITE/epitraxr"
